# Add "% of Q Drop's" column (I1) to the grade distribution header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "% of Q Drop's"
